# Fix up the descriptions on the "read.me" legend sheet (column B) and
# widen column B to fit the new, longer text.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("read.me")

$ws.Range("B2").Value = "raw colony count of glycerol stock into BHI broth after 72hrs, 37 C, shaken. ""bhi_6"" refers to 10^-6 dilution VRB media"
$ws.Range("B3").Value = "E. coli added from glycerol to BHI broth, then broth culutre to mik. Negative is nothing added to broth or milk."
$ws.Range("B4").Value = "raw colony count of E. coli culture taken from BHI into milk, fermented for 72hrs, 28 C, unshaken. ""milk_4"" refers to 10^-4 dilution on VRB media"

$ws.Columns.Item(2).ColumnWidth = 116.49869791666667

$ws.Activate()
$ws.Range("B3").Select()
